$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.895.01"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.230.31"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "2.555.64"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.860"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.220.26"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "41.700.50"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.76%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0723"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0297"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.50%  "
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.07%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.209"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
